$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values are purely numeric-looking
# strings, so Excel keeps them as literal text instead of converting to floats
# (matching the source data which stores all Price/Volume cells as text).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '60.258.50'
$ws.Range('E2').Value = '  -3.80%  '
$ws.Range('D3').Value = '3.304.58'
$ws.Range('E3').Value = '  -4.18%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '559.89'
$ws.Range('E5').Value = '  -3.18%  '
$ws.Range('D6').Value = '143.64'
$ws.Range('E6').Value = '  -3.33%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '3.305.46'
$ws.Range('E8').Value = '  -4.07%  '
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('D10').Value = '7.78'
$ws.Range('E10').Value = '  -3.17%  '
$ws.Range('E11').Value = '  -3.41%  '
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').Value = '3.877.57'
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('E15').Value = '  -3.89%  '
$ws.Range('D16').Value = '3.315.07'
$ws.Range('E16').Value = '  -3.84%  '
$ws.Range('E17').Value = '  -3.45%  '
$ws.Range('D18').Value = '60.314.73'
$ws.Range('E18').Value = '  -3.78%  '
$ws.Range('E19').Value = '  -3.95%  '
$ws.Range('D20').Value = '14.35'
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').Value = '8.61'
$ws.Range('E21').Value = '  -4.20%  '
$ws.Range('D22').Value = '373.36'
$ws.Range('E22').Value = '  -3.61%  '
$ws.Range('D23').Value = '74.15'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('E24').Value = '  -2.98%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = '3.469.89'
$ws.Range('E26').Value = '  -3.20%  '
$ws.Range('E27').Value = '  -8.31%  '
$ws.Range('E28').Value = '  -4.89%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '7.21'
$ws.Range('E30').Value = '  -5.99%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('D33').Value = '7.63'
$ws.Range('E33').Value = '  -4.61%  '
$ws.Range('E34').Value = '  -2.75%  '
$ws.Range('E35').Value = '  -3.98%  '
$ws.Range('D36').Value = '5.16'
$ws.Range('E36').Value = '  -3.97%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '166.05'
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = '6.76'
$ws.Range('E38').Value = '  -2.69%  '
$ws.Range('D39').Value = '1.52'
$ws.Range('E39').Value = '  -6.33%  '
$ws.Range('D40').Value = '27.77'
$ws.Range('E40').Value = '  -13.52%  '
$ws.Range('E42').Value = '  -4.90%  '
$ws.Range('E44').Value = '  -4.33%  '
$ws.Range('D45').Value = '4.21'
$ws.Range('E45').Value = '  -3.95%  '
$ws.Range('E46').Value = '  -4.85%  '
$ws.Range('E47').Value = '  -4.21%  '
$ws.Range('D48').Value = '2.375.98'
$ws.Range('E48').Value = '  -7.52%  '
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('D51').Value = '21.66'
$ws.Range('E51').Value = '  -3.84%  '
